$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.848.67'
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").Value = '1.635.99'
$ws.Range("E3").Value = '  -0.18%  '

$ws.Range("E4").Value = '  -0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.76'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.510'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +1.70%  '

$ws.Range("E7").Value = '  -0.51%  '

$ws.Range("E8").Value = '  +1.34%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0625'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +0.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.01'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +3.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.14%  '

$ws.Range("D12").Value = '1.865.25'
$ws.Range("E12").Value = '  -0.23%  '

$ws.Range("D13").Value = '1.640.04'
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("E14").Value = '  -0.73%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.529'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.60%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.58'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +2.30%  '

$ws.Range("D17").Value = '26.849.10'
$ws.Range("E17").Value = '  +0.09%  '

$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("E18").Value = '  -0.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '219.44'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +1.05%  '

$ws.Range("E20").Value = '  -0.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.77'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +2.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.39'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +0.70%  '

$ws.Range("E23").Value = '  +3.86%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.16'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.67'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -0.29%  '

$ws.Range("E26").Value = '  -0.28%  '

$ws.Range("E27").Value = '  +3.33%  '

$ws.Range("E28").Value = '  +0.36%  '

$ws.Range("E29").Value = '  +0.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0503'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -0.55%  '

$ws.Range("E31").Value = '  -1.49%  '

$ws.Range("E32").Value = '  -1.45%  '

$ws.Range("E33").Value = '  +0.28%  '

$ws.Range("E34").Value = '  +0.54%  '

$ws.Range("D35").Value = '1.253.21'
$ws.Range("E35").Value = '  -0.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.44'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -0.16%  '

$ws.Range("E37").Value = '  +1.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.534'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +0.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.831'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +1.67%  '

$ws.Range("E40").Value = '  -0.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.808'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +0.28%  '

$ws.Range("E42").Value = '  +0.89%  '

$ws.Range("D43").Value = '1.775.97'
$ws.Range("E43").Value = '  -0.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.66'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +1.73%  '

$ws.Range("E45").Value = '  -1.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.54'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -0.80%  '

$ws.Range("E47").Value = '  -0.63%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0105'
$ws.Range("E48").Value = '  +2.42%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0513'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -0.50%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.63'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +1.28%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0959'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -0.73%  '

